# Generate Report for Handback
# Adds a new handback record (6a72a293-1cef-46b6-89c7-c0b0778ecde2.md) as row 4
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$fileId   = "6a72a293-1cef-46b6-89c7-c0b0778ecde2"
$fileName = "$fileId.md"
$pathName = "e2e\$fileName"
$ext      = ".md"
$statusSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $fileName
$wsOverview.Range("B4").Value = $pathName
$wsOverview.Range("C4").Value = $ext
$wsOverview.Range("E4").Value = $statusSync
$wsOverview.Range("F4").Value = $statusSync
$wsOverview.Range("G4").Value = "2016-09-04 04:49:20"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4843fa84924739337db6b520ba98a5c9b878132/e2e/$fileName",
    "",
    "",
    $pathName
) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$zhXlf = "$fileId.f4843fa84924739337db6b520ba98a5c9b878132.zh-cn.xlf"

$wsZhCn.Range("A4").Value = $fileName
$wsZhCn.Range("B4").Value = $ext
$wsZhCn.Range("C4").Value = $statusSync
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "True"
$wsZhCn.Range("G4").Value = $zhXlf
$wsZhCn.Range("H4").Value = "2016-09-04 04:49:15"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I4").Value = $fileName
$wsZhCn.Range("J4").Value = $zhXlf
$wsZhCn.Range("K4").Value = "2016-09-04 04:49:32"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L4").Value = ""
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "False"
$wsZhCn.Range("P4").Value = ""

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4843fa84924739337db6b520ba98a5c9b878132/e2e/$fileName",
    "",
    "",
    $fileName
) | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f4843fa84924739337db6b520ba98a5c9b878132/e2e/$fileName",
    "",
    "",
    $fileName
) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$deXlf = "$fileId.f4843fa84924739337db6b520ba98a5c9b878132.de-de.xlf"

$wsDeDe.Range("A4").Value = $fileName
$wsDeDe.Range("B4").Value = $ext
$wsDeDe.Range("C4").Value = $statusSync
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "True"
$wsDeDe.Range("G4").Value = $deXlf
$wsDeDe.Range("H4").Value = "2016-09-04 04:49:20"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I4").Value = $fileName
$wsDeDe.Range("J4").Value = $deXlf
$wsDeDe.Range("K4").Value = "2016-09-04 04:49:39"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L4").Value = ""
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "False"
$wsDeDe.Range("P4").Value = ""

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4843fa84924739337db6b520ba98a5c9b878132/e2e/$fileName",
    "",
    "",
    $fileName
) | Out-Null
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f4843fa84924739337db6b520ba98a5c9b878132/e2e/$fileName",
    "",
    "",
    $fileName
) | Out-Null
